$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: duplicate the last existing triplet (rows 75:77, one full "game session" block)
# five times to create rows 78-92, preserving per-row relative formulas and cell styles
# (e.g. the date style on column B) exactly as a real copy/insert-rows operation would.
$ws.Rows("75:77").Copy()
$ws.Rows("78:80").Insert()
$ws.Rows("75:77").Copy()
$ws.Rows("81:83").Insert()
$ws.Rows("75:77").Copy()
$ws.Rows("84:86").Insert()
$ws.Rows("75:77").Copy()
$ws.Rows("87:89").Insert()
$ws.Rows("75:77").Copy()
$ws.Rows("90:92").Insert()

# Step 2: overwrite the copied cells with the actual new season data for each row.
# --- row 78 ---
$ws.Range("A78").Value = 2022
$ws.Range("B78").Value = 44924
$ws.Range("C78").Value = 26
$ws.Range("D78").Value = "Husum (Memeler Str.)"
$ws.Range("E78").Value = "MF"
$ws.Range("F78").Value = 11
$ws.Range("G78").Value = 4
$ws.Range("H78").Value = 2
$ws.Range("I78").Value = "first"
$ws.Range("J78").Value = "outer"
$ws.Range("L78").Value = 5
$ws.Range("M78").Value = 5
$ws.Range("N78").Value = 3
$ws.Range("O78").Value = 6
$ws.Range("P78").Value = 5
$ws.Range("Q78").Value = 3
$ws.Range("R78").Value = 6
$ws.Range("S78").Value = 3
$ws.Range("T78").Value = 2
$ws.Range("U78").Value = 0
$ws.Range("V78").Value = 0
$ws.Range("W78").Value = 0
$ws.Range("X78").Value = 0
$ws.Range("Y78").Value = 0
$ws.Range("Z78").Value = 0
$ws.Range("AA78").Value = 0
$ws.Range("AB78").Value = 0
$ws.Range("AC78").Value = 0
$ws.Range("AD78").Value = 0
$ws.Range("AE78").Value = 1
$ws.Range("AF78").Value = 0
$ws.Range("AG78").Value = 0

# --- row 79 ---
$ws.Range("A79").Value = 2022
$ws.Range("B79").Value = 44924
$ws.Range("C79").Value = 26
$ws.Range("D79").Value = "Husum (Memeler Str.)"
$ws.Range("E79").Value = "JHC"
$ws.Range("F79").Value = 13
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 1
$ws.Range("I79").Value = "second"
$ws.Range("J79").Value = "outer"
$ws.Range("L79").Value = 3
$ws.Range("M79").Value = 5
$ws.Range("N79").Value = 4
$ws.Range("O79").Value = 11
$ws.Range("P79").Value = 5
$ws.Range("Q79").Value = 6
$ws.Range("R79").Value = 5
$ws.Range("S79").Value = 0
$ws.Range("T79").Value = 5
$ws.Range("U79").Value = 0
$ws.Range("V79").Value = 0
$ws.Range("W79").Value = 0
$ws.Range("X79").Value = 0
$ws.Range("Y79").Value = 0
$ws.Range("Z79").Value = 0
$ws.Range("AA79").Value = 0
$ws.Range("AB79").Value = 0
$ws.Range("AC79").Value = 1
$ws.Range("AD79").Value = 0
$ws.Range("AE79").Value = 0
$ws.Range("AF79").Value = 0
$ws.Range("AG79").Value = 0

# --- row 80 ---
$ws.Range("A80").Value = 2022
$ws.Range("B80").Value = 44924
$ws.Range("C80").Value = 26
$ws.Range("D80").Value = "Husum (Memeler Str.)"
$ws.Range("E80").Value = "PF"
$ws.Range("F80").Value = 8
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 3
$ws.Range("I80").Value = "third"
$ws.Range("J80").Value = "outer"
$ws.Range("L80").Value = 6
$ws.Range("M80").Value = 5
$ws.Range("N80").Value = 2
$ws.Range("O80").Value = 4
$ws.Range("P80").Value = 2
$ws.Range("Q80").Value = 11
$ws.Range("R80").Value = 5
$ws.Range("S80").Value = 0
$ws.Range("T80").Value = 5
$ws.Range("U80").Value = 2
$ws.Range("V80").Value = 0
$ws.Range("W80").Value = 0
$ws.Range("X80").Value = 0
$ws.Range("Y80").Value = 0
$ws.Range("Z80").Value = 0
$ws.Range("AA80").Value = 0
$ws.Range("AB80").Value = 0
$ws.Range("AC80").Value = 0
$ws.Range("AD80").Value = 0
$ws.Range("AE80").Value = 0
$ws.Range("AF80").Value = 0
$ws.Range("AG80").Value = 0

# --- row 81 ---
$ws.Range("A81").Value = 2022
$ws.Range("B81").Value = 44924
$ws.Range("C81").Value = 27
$ws.Range("D81").Value = "Husum (Memeler Str.)"
$ws.Range("E81").Value = "PF"
$ws.Range("F81").Value = 8
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 2
$ws.Range("I81").Value = "first"
$ws.Range("J81").Value = "outer"
$ws.Range("L81:AG81").ClearContents()

# --- row 82 ---
$ws.Range("A82").Value = 2022
$ws.Range("B82").Value = 44924
$ws.Range("C82").Value = 27
$ws.Range("D82").Value = "Husum (Memeler Str.)"
$ws.Range("E82").Value = "MF"
$ws.Range("F82").Value = 5
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 3
$ws.Range("I82").Value = "second"
$ws.Range("J82").Value = "outer"
$ws.Range("L82:AG82").ClearContents()

# --- row 83 ---
$ws.Range("A83").Value = 2022
$ws.Range("B83").Value = 44924
$ws.Range("C83").Value = 27
$ws.Range("D83").Value = "Husum (Memeler Str.)"
$ws.Range("E83").Value = "JHC"
$ws.Range("F83").Value = 14
$ws.Range("G83").Value = 5
$ws.Range("H83").Value = 1
$ws.Range("I83").Value = "third"
$ws.Range("J83").Value = "outer"
$ws.Range("L83:AG83").ClearContents()

# --- row 84 ---
$ws.Range("A84").Value = 2022
$ws.Range("B84").Value = 44924
$ws.Range("C84").Value = 28
$ws.Range("D84").Value = "Husum (Memeler Str.)"
$ws.Range("E84").Value = "MF"
$ws.Range("F84").Value = 7
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 3
$ws.Range("I84").Value = "first"
$ws.Range("J84").Value = "outer"
$ws.Range("L84").Value = 5
$ws.Range("M84").Value = 5
$ws.Range("N84").Value = 5
$ws.Range("O84").Value = 3
$ws.Range("P84").Value = 3
$ws.Range("Q84").Value = 12
$ws.Range("R84").Value = 2
$ws.Range("S84").Value = 0
$ws.Range("T84").Value = 2
$ws.Range("U84").Value = 3
$ws.Range("V84").Value = 0
$ws.Range("W84").Value = 0
$ws.Range("X84").Value = 0
$ws.Range("Y84").Value = 0
$ws.Range("Z84").Value = 0
$ws.Range("AA84").Value = 0
$ws.Range("AB84").Value = 0
$ws.Range("AC84").Value = 0
$ws.Range("AD84").Value = 0
$ws.Range("AE84").Value = 1
$ws.Range("AF84").Value = 0
$ws.Range("AG84").Value = 0

# --- row 85 ---
$ws.Range("A85").Value = 2022
$ws.Range("B85").Value = 44924
$ws.Range("C85").Value = 28
$ws.Range("D85").Value = "Husum (Memeler Str.)"
$ws.Range("E85").Value = "PF"
$ws.Range("F85").Value = 13
$ws.Range("G85").Value = 6
$ws.Range("H85").Value = 1
$ws.Range("I85").Value = "second"
$ws.Range("J85").Value = "outer"
$ws.Range("L85").Value = 6
$ws.Range("M85").Value = 5
$ws.Range("N85").Value = 5
$ws.Range("O85").Value = 3
$ws.Range("P85").Value = 6
$ws.Range("Q85").Value = 4
$ws.Range("R85").Value = 4
$ws.Range("S85").Value = 5
$ws.Range("T85").Value = 4
$ws.Range("U85").Value = 1
$ws.Range("V85").Value = 0
$ws.Range("W85").Value = 0
$ws.Range("X85").Value = 0
$ws.Range("Y85").Value = 0
$ws.Range("Z85").Value = 0
$ws.Range("AA85").Value = 0
$ws.Range("AB85").Value = 0
$ws.Range("AC85").Value = 0
$ws.Range("AD85").Value = 1
$ws.Range("AE85").Value = 0
$ws.Range("AF85").Value = 0
$ws.Range("AG85").Value = 0

# --- row 86 ---
$ws.Range("A86").Value = 2022
$ws.Range("B86").Value = 44924
$ws.Range("C86").Value = 28
$ws.Range("D86").Value = "Husum (Memeler Str.)"
$ws.Range("E86").Value = "JHC"
$ws.Range("F86").Value = 10
$ws.Range("G86").Value = 4
$ws.Range("H86").Value = 2
$ws.Range("I86").Value = "third"
$ws.Range("J86").Value = "outer"
$ws.Range("L86").Value = 5
$ws.Range("M86").Value = 5
$ws.Range("N86").Value = 7
$ws.Range("O86").Value = 3
$ws.Range("P86").Value = 3
$ws.Range("Q86").Value = 8
$ws.Range("R86").Value = 3
$ws.Range("S86").Value = 5
$ws.Range("T86").Value = 0
$ws.Range("U86").Value = 3
$ws.Range("V86").Value = 0
$ws.Range("W86").Value = 0
$ws.Range("X86").Value = 0
$ws.Range("Y86").Value = 0
$ws.Range("Z86").Value = 0
$ws.Range("AA86").Value = 0
$ws.Range("AB86").Value = 0
$ws.Range("AC86").Value = 1
$ws.Range("AD86").Value = 0
$ws.Range("AE86").Value = 0
$ws.Range("AF86").Value = 0
$ws.Range("AG86").Value = 0

# --- row 87 ---
$ws.Range("A87").Value = 2022
$ws.Range("B87").Value = 44924
$ws.Range("C87").Value = 29
$ws.Range("D87").Value = "Husum (Memeler Str.)"
$ws.Range("E87").Value = "JHC"
$ws.Range("F87").Value = 8
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 3
$ws.Range("I87").Value = "first"
$ws.Range("J87").Value = "outer"
$ws.Range("L87").Value = 4
$ws.Range("M87").Value = 4
$ws.Range("N87").Value = 3
$ws.Range("O87").Value = 15
$ws.Range("P87").Value = 3
$ws.Range("Q87").Value = 0
$ws.Range("R87").Value = 4
$ws.Range("S87").Value = 3
$ws.Range("T87").Value = 0
$ws.Range("U87").Value = 0
$ws.Range("V87").Value = 0
$ws.Range("W87").Value = 0
$ws.Range("X87").Value = 0
$ws.Range("Y87").Value = 0
$ws.Range("Z87").Value = 0
$ws.Range("AA87").Value = 1
$ws.Range("AB87").Value = 0
$ws.Range("AC87").Value = 0
$ws.Range("AD87").Value = 0
$ws.Range("AE87").Value = 0
$ws.Range("AF87").Value = 1
$ws.Range("AG87").Value = 0

# --- row 88 ---
$ws.Range("A88").Value = 2022
$ws.Range("B88").Value = 44924
$ws.Range("C88").Value = 29
$ws.Range("D88").Value = "Husum (Memeler Str.)"
$ws.Range("E88").Value = "PF"
$ws.Range("F88").Value = 9
$ws.Range("G88").Value = 4
$ws.Range("H88").Value = 2
$ws.Range("I88").Value = "second"
$ws.Range("J88").Value = "outer"
$ws.Range("L88").Value = 6
$ws.Range("M88").Value = 4
$ws.Range("N88").Value = 6
$ws.Range("O88").Value = 5
$ws.Range("P88").Value = 0
$ws.Range("Q88").Value = 8
$ws.Range("R88").Value = 2
$ws.Range("S88").Value = 6
$ws.Range("T88").Value = 0
$ws.Range("U88").Value = 0
$ws.Range("V88").Value = 0
$ws.Range("W88").Value = 0
$ws.Range("X88").Value = 0
$ws.Range("Y88").Value = 0
$ws.Range("Z88").Value = 0
$ws.Range("AA88").Value = 0
$ws.Range("AB88").Value = 0
$ws.Range("AC88").Value = 0
$ws.Range("AD88").Value = 0
$ws.Range("AE88").Value = 0
$ws.Range("AF88").Value = 0
$ws.Range("AG88").Value = 1

# --- row 89 ---
$ws.Range("A89").Value = 2022
$ws.Range("B89").Value = 44924
$ws.Range("C89").Value = 29
$ws.Range("D89").Value = "Husum (Memeler Str.)"
$ws.Range("E89").Value = "MF"
$ws.Range("F89").Value = 13
$ws.Range("G89").Value = 7
$ws.Range("H89").Value = 1
$ws.Range("I89").Value = "third"
$ws.Range("J89").Value = "outer"
$ws.Range("L89").Value = 5
$ws.Range("M89").Value = 4
$ws.Range("N89").Value = 7
$ws.Range("O89").Value = 0
$ws.Range("P89").Value = 8
$ws.Range("Q89").Value = 4
$ws.Range("R89").Value = 4
$ws.Range("S89").Value = 7
$ws.Range("T89").Value = 4
$ws.Range("U89").Value = 0
$ws.Range("V89").Value = 0
$ws.Range("W89").Value = 0
$ws.Range("X89").Value = 0
$ws.Range("Y89").Value = 0
$ws.Range("Z89").Value = 0
$ws.Range("AA89").Value = 0
$ws.Range("AB89").Value = 0
$ws.Range("AC89").Value = 0
$ws.Range("AD89").Value = 1
$ws.Range("AE89").Value = 1
$ws.Range("AF89").Value = 0
$ws.Range("AG89").Value = 0

# --- row 90 ---
$ws.Range("A90").Value = 2022
$ws.Range("B90").Value = 44924
$ws.Range("C90").Value = 30
$ws.Range("D90").Value = "Husum (Memeler Str.)"
$ws.Range("E90").Value = "PF"
$ws.Range("F90").Value = 13
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 1
$ws.Range("I90").Value = "first"
$ws.Range("J90").Value = "outer"
$ws.Range("L90").Value = 6
$ws.Range("M90").Value = 5
$ws.Range("N90").Value = 5
$ws.Range("O90").Value = 4
$ws.Range("P90").Value = 5
$ws.Range("Q90").Value = 4
$ws.Range("R90").Value = 5
$ws.Range("S90").Value = 5
$ws.Range("T90").Value = 0
$ws.Range("U90").Value = 5
$ws.Range("V90").Value = 0
$ws.Range("W90").Value = 0
$ws.Range("X90").Value = 0
$ws.Range("Y90").Value = 0
$ws.Range("Z90").Value = 0
$ws.Range("AA90").Value = 0
$ws.Range("AB90").Value = 0
$ws.Range("AC90").Value = 0
$ws.Range("AD90").Value = 1
$ws.Range("AE90").Value = 0
$ws.Range("AF90").Value = 1
$ws.Range("AG90").Value = 0

# --- row 91 ---
$ws.Range("A91").Value = 2022
$ws.Range("B91").Value = 44924
$ws.Range("C91").Value = 30
$ws.Range("D91").Value = "Husum (Memeler Str.)"
$ws.Range("E91").Value = "JHC"
$ws.Range("F91").Value = 10
$ws.Range("G91").Value = 4
$ws.Range("H91").Value = 3
$ws.Range("I91").Value = "second"
$ws.Range("J91").Value = "outer"
$ws.Range("L91").Value = 6
$ws.Range("M91").Value = 5
$ws.Range("N91").Value = 3
$ws.Range("O91").Value = 8
$ws.Range("P91").Value = 8
$ws.Range("Q91").Value = 2
$ws.Range("R91").Value = 4
$ws.Range("S91").Value = 3
$ws.Range("T91").Value = 0
$ws.Range("U91").Value = 3
$ws.Range("V91").Value = 0
$ws.Range("W91").Value = 0
$ws.Range("X91").Value = 0
$ws.Range("Y91").Value = 0
$ws.Range("Z91").Value = 0
$ws.Range("AA91").Value = 0
$ws.Range("AB91").Value = 0
$ws.Range("AC91").Value = 1
$ws.Range("AD91").Value = 0
$ws.Range("AE91").Value = 0
$ws.Range("AF91").Value = 0
$ws.Range("AG91").Value = 0

# --- row 92 ---
$ws.Range("A92").Value = 2022
$ws.Range("B92").Value = 44924
$ws.Range("C92").Value = 30
$ws.Range("D92").Value = "Husum (Memeler Str.)"
$ws.Range("E92").Value = "MF"
$ws.Range("F92").Value = 11
$ws.Range("G92").Value = 4
$ws.Range("H92").Value = 2
$ws.Range("I92").Value = "third"
$ws.Range("J92").Value = "outer"
$ws.Range("L92").Value = 6
$ws.Range("M92").Value = 4
$ws.Range("N92").Value = 12
$ws.Range("O92").Value = 0
$ws.Range("P92").Value = 3
$ws.Range("Q92").Value = 2
$ws.Range("R92").Value = 4
$ws.Range("S92").Value = 12
$ws.Range("T92").Value = 0
$ws.Range("U92").Value = 0
$ws.Range("V92").Value = 0
$ws.Range("W92").Value = 0
$ws.Range("X92").Value = 0
$ws.Range("Y92").Value = 0
$ws.Range("Z92").Value = 0
$ws.Range("AA92").Value = 0
$ws.Range("AB92").Value = 0
$ws.Range("AC92").Value = 0
$ws.Range("AD92").Value = 0
$ws.Range("AE92").Value = 0
$ws.Range("AF92").Value = 0
$ws.Range("AG92").Value = 1

# Step 3: restore the on-screen selection to where the user ended up editing.
$ws.Range("H84").Select()
